$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.878.34"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.38%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.541.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.57%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.73"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.74%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.72%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.544"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.92"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0824"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.63"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.87%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.933.24"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.530.30"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.25"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +7.22%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.66%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.886.37"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.19"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.52%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.56%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.70"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.70%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.45"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.94"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.06"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.71"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.57%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.16%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +9.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.21"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.60"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.17"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.86"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.34%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0802"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.00"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.76%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.85%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.16%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.50"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.73%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.12"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.47"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.90"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.082.21"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.37%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.18%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.36"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.05"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.790.36"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.71"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.69%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.49%  "
